$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

# Summer 23 week 5 inputs - updated matchup average values
$ws.Range("D3").Value = 10.63
$ws.Range("E3").Value = 10.3

$ws.Range("C4").Value = 9.369999999999999
$ws.Range("F4").Value = 10

$ws.Range("C5").Value = 9.699999999999999
$ws.Range("F5").Value = 10.28
$ws.Range("G5").Value = 9.27
$ws.Range("H5").Value = 7.94

$ws.Range("D6").Value = 10
$ws.Range("E6").Value = 9.720000000000001

$ws.Range("E7").Value = 10.73
$ws.Range("J7").Value = 8.33

$ws.Range("E8").Value = 12.06

$ws.Range("G10").Value = 11.67
